$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.015.24'
$ws.Range("E2").Value = '  -2.29%  '
$ws.Range("D3").Value = '3.121.54'
$ws.Range("E3").Value = '  -0.80%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.89'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.30'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.14%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.119.30'
$ws.Range("E8").Value = '  -0.74%  '
$ws.Range("E9").Value = '  -1.69%  '
$ws.Range("E10").Value = '  -3.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.24'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.455'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000245'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.08'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.65%  '
$ws.Range("D15").Value = '3.634.48'
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("E16").Value = '  +1.69%  '
$ws.Range("D17").Value = '63.020.85'
$ws.Range("E17").Value = '  -2.26%  '
$ws.Range("D18").Value = '3.125.27'
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.65'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '471.17'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.89%  '
$ws.Range("E21").Value = '  -3.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.695'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.65'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.71'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.91'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.77%  '
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("E27").Value = '  -1.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.86'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.83'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.00%  '
$ws.Range("E30").Value = '  +1.56%  '
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.67'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.107'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.94%  '
$ws.Range("E34").Value = '  -5.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.07'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.78'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.91'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.92%  '
$ws.Range("D38").Value = '0.0₃0693'
$ws.Range("E38").Value = '  -10.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0387'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '417.37'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.38%  '
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").Value = '2.896.91'
$ws.Range("E42").Value = '  +1.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.65'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -12.56%  '
$ws.Range("E44").Value = '  -5.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.263'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.51%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.10'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.37'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.112'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.45%  '
$ws.Range("E50").Value = '  -7.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.96'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.00%  '
